$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 04:54:55"
$wsZhCn.Range("H2").Value = "2016-03-20 04:55:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 04:54:59"
$wsDeDe.Range("H2").Value = "2016-03-20 04:55:21"
